$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryList")

# --- Row 4 ---
$ws.Range("C4").Value = "In0001"
$ws.Range("D4").Value = "Tovar"
$ws.Range("E4").Value = "Item 1"
$ws.Range("F4").Value = 33
$ws.Range("I4").Value = 200
$ws.Range("M4").Value = "Item 2"

# --- Row 5 ---
$ws.Range("C5").Value = "In0002"
$ws.Range("D5").Value = "Tovar1"
$ws.Range("E5").Value = "Item 2"
$ws.Range("F5").Value = 34
$ws.Range("I5").Value = 20
$ws.Range("M5").Value = "Item 1"

# --- Row 6 ---
# D6, E6 and I6 already exist as blank cells with no explicit style, so a
# plain .Value assignment would leave them with the default (unstyled)
# format. Copy the surrounding column formatting first (format-paste only)
# so the underlying style index lines up with the rest of the column,
# then set the values.
$ws.Range("M4").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("L6").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("C6").Value = "In0003"
$ws.Range("D6").Value = "Tovar2"
$ws.Range("E6").Value = "Item 1"
$ws.Range("F6").Value = 123
$ws.Range("I6").Value = 132
$ws.Range("M6").Value = "Item 1"
